$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfrsf14"
$ws.Range("C2").Value = "Cd160"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.118296000000001
$ws.Range("H2").Value = 21.354888
$ws.Range("I2").Value = 0.4484208155031491
$ws.Range("J2").Value = 0.4484208155031491
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 1.093326
$ws.Range("N2").Value = 3.279978
$ws.Range("O2").Value = 0.397438973321003
$ws.Range("P2").Value = 0.397438973321003
$ws.Range("Q2").Value = 7.782618092496001
$ws.Range("R2").Value = 70.043562832464
$ws.Range("S2").Value = 0.1782199085293385
$ws.Range("T2").Value = 0.1782199085293385

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfrsf14"
$ws.Range("C3").Value = "Cd160"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.118296000000001
$ws.Range("H3").Value = 21.354888
$ws.Range("I3").Value = 0.4484208155031491
$ws.Range("J3").Value = 0.4484208155031491
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.386846666666666
$ws.Range("N3").Value = 4.160539999999999
$ws.Range("O3").Value = 0.5041377552108476
$ws.Range("P3").Value = 0.5041377552108477
$ws.Range("Q3").Value = 9.871985079946665
$ws.Range("R3").Value = 88.84786571951999
$ws.Range("S3").Value = 0.2260658633175752
$ws.Range("T3").Value = 0.2260658633175752

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tnfrsf14"
$ws.Range("C4").Value = "Cd160"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.118296000000001
$ws.Range("H4").Value = 21.354888
$ws.Range("I4").Value = 0.4484208155031491
$ws.Range("J4").Value = 0.4484208155031491
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2707553333333333
$ws.Range("N4").Value = 0.8122659999999999
$ws.Range("O4").Value = 0.09842327146814941
$ws.Range("P4").Value = 0.09842327146814943
$ws.Range("Q4").Value = 1.927316606245333
$ws.Range("R4").Value = 17.345849456208
$ws.Range("S4").Value = 0.04413504365623539
$ws.Range("T4").Value = 0.04413504365623539

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfrsf14"
$ws.Range("C5").Value = "Cd160"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.371014666666666
$ws.Range("H5").Value = 19.113044
$ws.Range("I5").Value = 0.4013454332903815
$ws.Range("J5").Value = 0.4013454332903815
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 1.093326
$ws.Range("N5").Value = 3.279978
$ws.Range("O5").Value = 0.397438973321003
$ws.Range("P5").Value = 0.397438973321003
$ws.Range("Q5").Value = 6.965595981448
$ws.Range("R5").Value = 62.690363833032
$ws.Range("S5").Value = 0.1595103169540023
$ws.Range("T5").Value = 0.1595103169540023

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tnfrsf14"
$ws.Range("C6").Value = "Cd160"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.371014666666666
$ws.Range("H6").Value = 19.113044
$ws.Range("I6").Value = 0.4013454332903815
$ws.Range("J6").Value = 0.4013454332903815
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.386846666666666
$ws.Range("N6").Value = 4.160539999999999
$ws.Range("O6").Value = 0.5041377552108476
$ws.Range("P6").Value = 0.5041377552108477
$ws.Range("Q6").Value = 8.835620453751108
$ws.Range("R6").Value = 79.52058408375999
$ws.Range("S6").Value = 0.2023333858031379
$ws.Range("T6").Value = 0.202333385803138

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tnfrsf14"
$ws.Range("C7").Value = "Cd160"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.371014666666666
$ws.Range("H7").Value = 19.113044
$ws.Range("I7").Value = 0.4013454332903815
$ws.Range("J7").Value = 0.4013454332903815
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2707553333333333
$ws.Range("N7").Value = 0.8122659999999999
$ws.Range("O7").Value = 0.09842327146814941
$ws.Range("P7").Value = 0.09842327146814943
$ws.Range("Q7").Value = 1.724986199744889
$ws.Range("R7").Value = 15.524875797704
$ws.Range("S7").Value = 0.03950173053324127
$ws.Range("T7").Value = 0.03950173053324128

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tnfrsf14"
$ws.Range("C8").Value = "Cd160"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.384832
$ws.Range("H8").Value = 7.154496
$ws.Range("I8").Value = 0.1502337512064693
$ws.Range("J8").Value = 0.1502337512064693
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 1.093326
$ws.Range("N8").Value = 3.279978
$ws.Range("O8").Value = 0.397438973321003
$ws.Range("P8").Value = 0.397438973321003
$ws.Range("Q8").Value = 2.607398831232
$ws.Range("R8").Value = 23.466589481088
$ws.Range("S8").Value = 0.05970874783766218
$ws.Range("T8").Value = 0.05970874783766218

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tnfrsf14"
$ws.Range("C9").Value = "Cd160"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.384832
$ws.Range("H9").Value = 7.154496
$ws.Range("I9").Value = 0.1502337512064693
$ws.Range("J9").Value = 0.1502337512064693
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.386846666666666
$ws.Range("N9").Value = 4.160539999999999
$ws.Range("O9").Value = 0.5041377552108476
$ws.Range("P9").Value = 0.5041377552108477
$ws.Range("Q9").Value = 3.307396309759999
$ws.Range("R9").Value = 29.76656678784
$ws.Range("S9").Value = 0.07573850609013442
$ws.Range("T9").Value = 0.07573850609013444

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tnfrsf14"
$ws.Range("C10").Value = "Cd160"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.384832
$ws.Range("H10").Value = 7.154496
$ws.Range("I10").Value = 0.1502337512064693
$ws.Range("J10").Value = 0.1502337512064693
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2707553333333333
$ws.Range("N10").Value = 0.8122659999999999
$ws.Range("O10").Value = 0.09842327146814941
$ws.Range("P10").Value = 0.09842327146814943
$ws.Range("Q10").Value = 0.6457059831039998
$ws.Range("R10").Value = 5.811353847935999
$ws.Range("S10").Value = 0.01478649727867275
$ws.Range("T10").Value = 0.01478649727867275
